$d = $word.ActiveDocument

function Edit-Sentence($findPhrase, $nounNew, $verbNew) {
    # Locate and clear the "<noun> <verb>" phrase (e.g. "respondent may" / "applicant may")
    $rng = $d.Content
    $found = $rng.Find.Execute($findPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Phrase not found: $findPhrase"
    }
    $rng.Text = ""
    $p = $rng.Start

    # Build: "<nounNew>" + " " + "<verbNew>" as three separate runs.
    $r1 = $d.Range($p, $p)
    $r1.InsertAfter($nounNew)
    $r1.Bold = 1

    $r2 = $d.Range($r1.End, $r1.End)
    $r2.InsertAfter(" ")
    $r2.Bold = 1

    $r3 = $d.Range($r2.End, $r2.End)
    $r3.InsertAfter($verbNew)
    $r3.Bold = 1

    # Remove the temporary bold formatting (reverse order) so the run boundaries
    # created above are preserved while the visible formatting matches the original.
    $r3.Bold = 0
    $r2.Bold = 0
    $r1.Bold = 0
}

Edit-Sentence "respondent may" "defendant" "should"
Edit-Sentence "applicant may" "claimant" "should"
